# Auto-generated from the unified OOXML diff for Typhon_Profits (8 worksheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Every write is a literal value assignment matching the target <v> in the diff; two cells are fully
# cleared (their <c> element disappears entirely) and two cells are newly populated (a <c> element appears
# that did not exist before).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 590.5
$ws.Range("I15").Value = 590.5
$ws.Range("K15").Value = 1771.5
$ws.Range("M15").Value = -1602.5
$ws.Range("H92").Value = 55555956
$ws.Range("I92").Value = 90909460
$ws.Range("K92").Value = 90909460
$ws.Range("M92").Value = -90908212
$ws.Range("H116").Value = 4651.0625
$ws.Range("I116").Value = 2126.25
$ws.Range("K116").Value = 2126.25
$ws.Range("M116").Value = 1315.75
$ws.Range("H129").Value = 233723.44
$ws.Range("J129").Value = 233723.44
$ws.Range("L129").Value = 701170.3200000001
$ws.Range("N129").Value = -711170.3200000001
$ws.Range("H132").Value = 3018.3438
$ws.Range("I132").Value = 3113.6072
$ws.Range("K132").Value = 9340.821599999999
$ws.Range("M132").Value = -6810.821599999999
$ws.Range("H137").Value = 1666.6
$ws.Range("I137").Value = 1400.2222
$ws.Range("K137").Value = 4200.6666
$ws.Range("M137").Value = -1650.6666
$ws.Range("H138").Value = 10103229
$ws.Range("J138").Value = 3330.5454
$ws.Range("L138").Value = 9991.636200000001
$ws.Range("N138").Value = -20271.6362
$ws.Range("H141").Value = 1101.3489
$ws.Range("I141").Value = 796.2564
$ws.Range("K141").Value = 2388.7692
$ws.Range("M141").Value = 2791.2308

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1984.22
$ws.Range("I32").Value = 1945.75
$ws.Range("J32").Value = 2426.625
$ws.Range("K32").Value = 1945.75
$ws.Range("L32").Value = 2426.625
$ws.Range("M32").Value = -1658.75
$ws.Range("N32").Value = -3000.625
$ws.Range("H41").Value = 2837.3333
$ws.Range("I41").Value = 2837.3333
$ws.Range("K41").Value = 2837.3333
$ws.Range("M41").Value = -2423.3333
$ws.Range("H63").Value = 2735
$ws.Range("J63").Value = 4000
$ws.Range("L63").Value = 4000
$ws.Range("N63").Value = -5372
$ws.Range("H66").Value = 2735
$ws.Range("J66").Value = 4000
$ws.Range("L66").Value = 20000
$ws.Range("N66").Value = -26864
$ws.Range("H74").Value = 37039252
$ws.Range("I74").Value = 37039252
$ws.Range("K74").Value = 37039252
$ws.Range("M74").Value = -37038378
$ws.Range("H77").Value = 37039252
$ws.Range("I77").Value = 37039252
$ws.Range("K77").Value = 185196260
$ws.Range("M77").Value = -185191892
$ws.Range("H122").Value = 2014.4615
$ws.Range("I122").Value = 1493.5238
$ws.Range("K122").Value = 4480.5714
$ws.Range("M122").Value = -2030.5714
$ws.Range("H132").Value = 12281.392
$ws.Range("I132").Value = 1367.0769
$ws.Range("K132").Value = 4101.2307
$ws.Range("M132").Value = -1571.2307

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H86").Value = 1552.6666
$ws.Range("I86").Value = 1391.5217
$ws.Range("K86").Value = 1391.5217
$ws.Range("M86").Value = -268.5217
$ws.Range("H89").Value = 1552.6666
$ws.Range("I89").Value = 1391.5217
$ws.Range("K89").Value = 6957.6085
$ws.Range("M89").Value = -1341.6085
$ws.Range("H134").Value = 3585.641
$ws.Range("I134").Value = 3951.3333
$ws.Range("J134").Value = 2366.6667
$ws.Range("K134").Value = 11853.9999
$ws.Range("L134").Value = 7100.000100000001
$ws.Range("M134").Value = -9318.999899999999
$ws.Range("N134").Value = -12170.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 628
$ws.Range("I16").Value = 487.375
$ws.Range("J16").Value = 853
$ws.Range("K16").Value = 487.375
$ws.Range("L16").Value = 853
$ws.Range("M16").Value = -200.375
$ws.Range("N16").Value = -1427
$ws.Range("H58").Value = 13978.474
$ws.Range("I58").Value = 828.4706
$ws.Range("K58").Value = 828.4706
$ws.Range("M58").Value = -625.4706
$ws.Range("H62").Value = 38464976
$ws.Range("I62").Value = 47622070
$ws.Range("J62").Value = 5162
$ws.Range("K62").Value = 47622070
$ws.Range("L62").Value = 5162
$ws.Range("M62").Value = -47621446
$ws.Range("N62").Value = -6410
$ws.Range("H65").Value = 38464976
$ws.Range("I65").Value = 47622070
$ws.Range("J65").Value = 5162
$ws.Range("K65").Value = 238110350
$ws.Range("L65").Value = 25810
$ws.Range("M65").Value = -238107230
$ws.Range("N65").Value = -32050
$ws.Range("H113").Value = 628
$ws.Range("I113").Value = 487.375
$ws.Range("J113").Value = 853
$ws.Range("K113").Value = 487.375
$ws.Range("L113").Value = 853
$ws.Range("M113").Value = 1682.625
$ws.Range("N113").Value = -5193
$ws.Range("H132").Value = 2189.111
$ws.Range("I132").Value = 1670.4117
$ws.Range("J132").Value = 11007
$ws.Range("K132").Value = 5011.2351
$ws.Range("L132").Value = 33021
$ws.Range("M132").Value = -2481.2351
$ws.Range("N132").Value = -38081
$ws.Range("H134").Value = 923.95746
$ws.Range("I134").Value = 783.60974
$ws.Range("K134").Value = 2350.82922
$ws.Range("M134").Value = 184.1707799999999
$ws.Range("H136").Value = 13978.474
$ws.Range("I136").Value = 828.4706
$ws.Range("K136").Value = 2485.4118
$ws.Range("M136").Value = 64.58820000000014

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 696.25
$ws.Range("J131").Value = 716.75824
$ws.Range("L131").Value = 2150.27472
$ws.Range("N131").Value = -12230.27472

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H70").Value = 11653.958
$ws.Range("I70").Value = 10288.929
$ws.Range("J70").Value = 13565
$ws.Range("K70").Value = 10288.929
$ws.Range("L70").Value = 13565
$ws.Range("M70").Value = -10018.929
$ws.Range("N70").Value = -14105
$ws.Range("H73").Value = 11653.958
$ws.Range("I73").Value = 10288.929
$ws.Range("J73").Value = 13565
$ws.Range("K73").Value = 10288.929
$ws.Range("L73").Value = 13565
$ws.Range("M73").Value = -9352.929
$ws.Range("N73").Value = -15437
$ws.Range("H80").Value = 3989.0527
$ws.Range("I80").Value = 3473.375
$ws.Range("J80").Value = 4364.091
$ws.Range("K80").Value = 3473.375
$ws.Range("L80").Value = 4364.091
$ws.Range("M80").Value = -2475.375
$ws.Range("N80").Value = -6360.091
$ws.Range("H83").Value = 3989.0527
$ws.Range("I83").Value = 3473.375
$ws.Range("J83").Value = 4364.091
$ws.Range("K83").Value = 17366.875
$ws.Range("L83").Value = 21820.455
$ws.Range("M83").Value = -12374.875
$ws.Range("N83").Value = -31804.455
$ws.Range("H102").Value = 29415160
$ws.Range("I102").Value = 29415160
$ws.Range("K102").Value = 29415160
$ws.Range("M102").Value = -29413538
$ws.Range("H132").Value = 93115.164
$ws.Range("I132").Value = 11738.2
$ws.Range("K132").Value = 35214.60000000001
$ws.Range("M132").Value = -32684.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 10000
$ws.Range("J9").Value = 10000
$ws.Range("L9").Value = 10000
$ws.Range("N9").Value = -10448
$ws.Range("H107").Value = 540
$ws.Range("I107").Value = 540
$ws.Range("K107").Value = 540
$ws.Range("M107").Value = 1380
$ws.Range("H132").Value = 309828.16
$ws.Range("I132").Value = 317960.47
$ws.Range("K132").Value = 953881.4099999999
$ws.Range("M132").Value = -951351.4099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4386
$ws.Range("J62").Value = 4540
$ws.Range("L62").Value = 4540
$ws.Range("N62").Value = -5788
$ws.Range("H65").Value = 4386
$ws.Range("J65").Value = 4540
$ws.Range("L65").Value = 22700
$ws.Range("N65").Value = -28940
$ws.Range("H81").Value = 66667800
$ws.Range("I81").Value = 1154.7693
$ws.Range("J81").Value = 500001000
$ws.Range("K81").Value = 2309.5386
$ws.Range("L81").Value = 1000002000
$ws.Range("M81").Value = -1248.5386
$ws.Range("N81").Value = -1000004122
$ws.Range("H84").Value = 66667800
$ws.Range("I84").Value = 1154.7693
$ws.Range("J84").Value = 500001000
$ws.Range("K84").Value = 11547.693
$ws.Range("L84").Value = 5000010000
$ws.Range("M84").Value = -6243.692999999999
$ws.Range("N84").Value = -5000020608
$ws.Range("H122").Value = 1283.52
$ws.Range("I122").Value = 1295.125
$ws.Range("K122").Value = 3885.375
$ws.Range("M122").Value = -1435.375
$ws.Range("H132").Value = 423.80283
$ws.Range("I132").Value = 423.125
$ws.Range("K132").Value = 1269.375
$ws.Range("M132").Value = 1260.625
$ws.Range("H136").Value = 16651242
$ws.Range("I136").Value = 23461316
$ws.Range("J136").Value = 4393.5
$ws.Range("K136").Value = 70383948
$ws.Range("L136").Value = 13180.5
$ws.Range("M136").Value = -70381398
$ws.Range("N136").Value = -18280.5
